# Remove the three decorative "solid fill" icon pictures (Inventory,
# Robot, Artificial Intelligence) that were added to the
# "Introduction to Reinforcement Learning" slide and then deleted again
# in the same authoring session.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

$namesToRemove = @("Graphic 10", "Graphic 9", "Graphic 12")

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($namesToRemove -contains $shp.Name) {
        $shp.Delete()
    }
}
